# Using CAPS for TOMES pattern syntax now.
#
# The "pattern" column (B) example text used by rows 4 and 5 of the
# "Entities" sheet switches from the lower-case "tomes_pattern: ..."
# syntax to the upper-case "TOMES_PATTERN: ..." syntax.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newPattern = 'TOMES_PATTERN: {"A","B"}, {"1","2"}'
$ws.Range("B4").Value = $newPattern
$ws.Range("B5").Value = $newPattern

# Re-select column B (the sheet was left with the whole "pattern" column
# selected: B1:B1048576, active cell B1).
$ws.Range("B1:B1048576").Select() | Out-Null

# Nudge the saved workbook window position down slightly, matching the
# recorded bookViews/workbookView yWindow value.
$excel.ActiveWindow.Top = 8550
